$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts N,O,P -> O,P,Q)
$ws.Columns("N").EntireColumn.Insert()

# The newly inserted column should inherit column M's width (11 chars, no bestFit)
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active/selected sheet
$ws.Activate()
$ws.Range("R6").Select()
